# Text updates as supplied by PM&C.
# Applies updated benchmark / description text to the "Description" sheet
# of the Legal Assistance Total Services workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# ---------------------------------------------------------------------
# Row 1 (B1): Benchmark - replace short benchmark text with the longer,
# more detailed benchmark description supplied by PM&C.
# ---------------------------------------------------------------------
$ws.Range("B1").Value2 = "95 per cent or more of representation services for legal aid commissions are delivered to people experiencing financial disadvantage " + [char]0x2013 + " to be achieved by each State in each six month period from 1 July 2015 onwards. 85 per cent or more of total representation services for community legal centres are delivered to people experiencing financial disadvantage " + [char]0x2013 + " to be achieved by each State in aggregate across all community legal centres in each six month period between 1 July 2015 and 30 June 2017."
$ws.Rows.Item(1).RowHeight = 85.05

# ---------------------------------------------------------------------
# Row 3: height adjustment only
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 13.8

# ---------------------------------------------------------------------
# Row 5 (B5): now holds a short sub-heading "Benchmark for legal aid
# commissions" instead of the old LAC paragraph.
# ---------------------------------------------------------------------
$ws.Range("B5").Value2 = "Benchmark for legal aid commissions"
$ws.Rows.Item(5).RowHeight = 13.8

# ---------------------------------------------------------------------
# Row 6 (B6): re-worded LAC paragraph (now sits under the new heading).
# ---------------------------------------------------------------------
$ws.Range("B6").Value2 = "All state and territory LACs met this benchmark for both reporting periods in 2015-16. Figures show that nationally, LACs are providing a consistently high proportion of representation services to financially disadvantaged people."
$ws.Rows.Item(6).RowHeight = 37.45

# ---------------------------------------------------------------------
# Row 7: drop the old empty A7 styled cell, turn B7 into the new
# "Benchmark for community legal centres" sub-heading (bold, italic,
# 12pt, dark colour, wrapped).
# ---------------------------------------------------------------------
$ws.Range("A7").Clear()

$ws.Range("B7").Value2 = "Benchmark for community legal centres"
$f7 = $ws.Range("B7").Font
$f7.Name = "Arial"
$f7.Size = 12
$f7.Bold = $true
$f7.Italic = $true
$f7.Color = 655360
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 15

# ---------------------------------------------------------------------
# Row 8 (new): CLC paragraph, re-worded.
# ---------------------------------------------------------------------
$ws.Range("B8").Value2 = "All states and territories met the benchmark for CLCs in 2015-16. Figures show that nationally, CLCs are providing a consistently high proportion of representation services to financially disadvantaged people. "
$f8 = $ws.Range("B8").Font
$f8.Name = "Arial"
$f8.Size = 12
$f8.Color = 655360
$ws.Range("B8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 39.7

# ---------------------------------------------------------------------
# Row 9 (new): forward-looking CLC benchmark change, re-worded
# ("90 per cent" instead of "90%").
# ---------------------------------------------------------------------
$ws.Range("B9").Value2 = "From 1 July 2017, the benchmark for CLCs will be increased so that 90 per cent or more of representation services are delivered to people experiencing financial disadvantage. "
$f9 = $ws.Range("B9").Font
$f9.Name = "Arial"
$f9.Size = 12
$f9.Color = 655360
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 26.95

# ---------------------------------------------------------------------
# Row 10 (new): Source row.
# ---------------------------------------------------------------------
$ws.Range("A10").Value2 = "Source"

$ws.Range("B10").Value2 = "National Partnership on Legal Assistance Services (NPLAS) reports provided by LACs."
$f10 = $ws.Range("B10").Font
$f10.Name = "Arial"
$f10.Size = 12
$f10.Color = 655360
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 15

$wb.Save()
